$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row (original data goes through row 20)
$lastRow = $ws.UsedRange.Rows.Count

# Capture the existing column-A labels (row 2..lastRow) before shifting columns,
# since these need to move into the new column B.
$labels = @{}
for ($r = 2; $r -le $lastRow; $r++) {
    $labels[$r] = $ws.Cells.Item($r, 1).Text
}

# Insert a new blank column before column B -> old B/C/D shift to C/D/E.
$ws.Range("B1").EntireColumn.Insert()

# New column B header, matching the bold/bordered header style used by the
# other header cells (copy format from the neighboring header cell).
$ws.Cells.Item(1, 2).Value = "segments"
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill column A with the numeric segment index (0-based) and column B with the
# label that used to live in column A.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
    $ws.Cells.Item($r, 2).Value = $labels[$r]
}
